$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 136.70967
$ws.Range("I33").Value = 90.92308
$ws.Range("K33").Value = 90.92308
$ws.Range("M33").Value = 138.07692
$ws.Range("H99").Value = 352.75
$ws.Range("I99").Value = 288.85715
$ws.Range("J99").Value = 800
$ws.Range("K99").Value = 866.5714499999999
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = 631.4285500000001
$ws.Range("N99").Value = -5396
$ws.Range("H104").Value = 662.8
$ws.Range("I104").Value = 662.8
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 1988.4
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -241.3999999999999
$ws.Range("N104").ClearContents()
$ws.Range("H111").Value = 2754.1
$ws.Range("I111").Value = 2674.4546
$ws.Range("J111").Value = 2851.4443
$ws.Range("K111").Value = 8023.3638
$ws.Range("L111").Value = 8554.332900000001
$ws.Range("M111").Value = -4956.3638
$ws.Range("N111").Value = -14688.3329
$ws.Range("H116").Value = 2516.6667
$ws.Range("I116").Value = 1975
$ws.Range("J116").Value = 3600
$ws.Range("K116").Value = 1975
$ws.Range("L116").Value = 3600
$ws.Range("M116").Value = 1467
$ws.Range("N116").Value = -10484
$ws.Range("H137").Value = 1825.2559
$ws.Range("I137").Value = 1694.7097
$ws.Range("J137").Value = 2162.5
$ws.Range("K137").Value = 5084.1291
$ws.Range("L137").Value = 6487.5
$ws.Range("M137").Value = -2534.1291
$ws.Range("N137").Value = -11587.5
$ws.Range("H138").Value = 4979.533
$ws.Range("I138").Value = 1987.579
$ws.Range("K138").Value = 5962.737
$ws.Range("M138").Value = -822.7370000000001
$ws.Range("H140").Value = 87682
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 87682
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 87682
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -98042

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H104").Value = 60000
$ws.Range("J104").Value = 60000
$ws.Range("L104").Value = 60000
$ws.Range("N104").Value = -66988
$ws.Range("H137").Value = 95000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 95000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 95000
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -105200

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2447.5588
$ws.Range("I31").Value = 1528.0555
$ws.Range("J31").Value = 3482
$ws.Range("K31").Value = 1528.0555
$ws.Range("L31").Value = 3482
$ws.Range("M31").Value = -1233.0555
$ws.Range("N31").Value = -4072
$ws.Range("H34").Value = 2447.5588
$ws.Range("I34").Value = 1528.0555
$ws.Range("J34").Value = 3482
$ws.Range("K34").Value = 1528.0555
$ws.Range("L34").Value = 3482
$ws.Range("M34").Value = -1326.0555
$ws.Range("N34").Value = -3886
$ws.Range("H94").Value = 869.2857
$ws.Range("I94").Value = 752
$ws.Range("J94").Value = 896.8823
$ws.Range("K94").Value = 752
$ws.Range("L94").Value = 896.8823
$ws.Range("M94").Value = -301
$ws.Range("N94").Value = -1798.8823
$ws.Range("H132").Value = 2385.6
$ws.Range("I132").Value = 2263.742
$ws.Range("J132").Value = 2655.4285
$ws.Range("K132").Value = 6791.226000000001
$ws.Range("L132").Value = 7966.2855
$ws.Range("M132").Value = -4261.226000000001
$ws.Range("N132").Value = -13026.2855
$ws.Range("H134").Value = 2308.1462
$ws.Range("I134").Value = 2096.1516
$ws.Range("K134").Value = 6288.4548
$ws.Range("M134").Value = -3753.4548

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 23824296
$ws.Range("I5").Value = 464
$ws.Range("J5").Value = 41692170
$ws.Range("K5").Value = 1392
$ws.Range("L5").Value = 125076510
$ws.Range("M5").Value = -1280
$ws.Range("N5").Value = -125076734
$ws.Range("H12").Value = 296.45
$ws.Range("I12").Value = 249.11111
$ws.Range("K12").Value = 747.3333299999999
$ws.Range("M12").Value = -574.3333299999999
$ws.Range("H33").Value = 78
$ws.Range("H107").Value = 978.2308
$ws.Range("J107").Value = 1167.7
$ws.Range("L107").Value = 3503.1
$ws.Range("N107").Value = -7343.1
$ws.Range("H132").Value = 1700
$ws.Range("I132").Value = 1300
$ws.Range("J132").Value = 1860
$ws.Range("K132").Value = 11700
$ws.Range("L132").Value = 16740
$ws.Range("M132").Value = -9170
$ws.Range("N132").Value = -21800
$ws.Range("H135").Value = 23824296
$ws.Range("I135").Value = 464
$ws.Range("J135").Value = 41692170
$ws.Range("K135").Value = 4176
$ws.Range("L135").Value = 375229530
$ws.Range("M135").Value = -1641
$ws.Range("N135").Value = -375234600

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8013.3335
$ws.Range("I80").Value = 12240
$ws.Range("J80").Value = 5900
$ws.Range("K80").Value = 12240
$ws.Range("L80").Value = 5900
$ws.Range("M80").Value = -11242
$ws.Range("N80").Value = -7896
$ws.Range("H83").Value = 8013.3335
$ws.Range("I83").Value = 12240
$ws.Range("J83").Value = 5900
$ws.Range("K83").Value = 61200
$ws.Range("L83").Value = 29500
$ws.Range("M83").Value = -56208
$ws.Range("N83").Value = -39484
$ws.Range("H113").Value = 3339.6
$ws.Range("I113").Value = 4332.6665
$ws.Range("K113").Value = 4332.6665
$ws.Range("M113").Value = -2162.6665
$ws.Range("H132").Value = 8798.333000000001
$ws.Range("I132").Value = 3427.4614
$ws.Range("K132").Value = 10282.3842
$ws.Range("M132").Value = -7752.3842

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 365.82352
$ws.Range("I55").Value = 314.875
$ws.Range("J55").Value = 411.1111
$ws.Range("K55").Value = 314.875
$ws.Range("L55").Value = 411.1111
$ws.Range("M55").Value = -141.875
$ws.Range("N55").Value = -757.1111000000001
$ws.Range("H61").Value = 1276088.6
$ws.Range("I61").Value = 40840.8
$ws.Range("K61").Value = 40840.8
$ws.Range("M61").Value = -40638.8
$ws.Range("H82").Value = 1229.1666
$ws.Range("J82").Value = 1361.1111
$ws.Range("L82").Value = 1361.1111
$ws.Range("N82").Value = -2083.1111
$ws.Range("H85").Value = 1229.1666
$ws.Range("J85").Value = 1361.1111
$ws.Range("L85").Value = 1361.1111
$ws.Range("N85").Value = -3857.1111
$ws.Range("H113").Value = 1276088.6
$ws.Range("I113").Value = 40840.8
$ws.Range("K113").Value = 40840.8
$ws.Range("M113").Value = -38670.8
$ws.Range("H122").Value = 7479.077
$ws.Range("I122").Value = 7229.304
$ws.Range("J122").Value = 7838.125
$ws.Range("K122").Value = 21687.912
$ws.Range("L122").Value = 23514.375
$ws.Range("M122").Value = -19237.912
$ws.Range("N122").Value = -28414.375
$ws.Range("H132").Value = 3040.05
$ws.Range("I132").Value = 2614.8064
$ws.Range("J132").Value = 4504.778
$ws.Range("K132").Value = 7844.4192
$ws.Range("L132").Value = 13514.334
$ws.Range("M132").Value = -5314.4192
$ws.Range("N132").Value = -18574.334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 32126
$ws.Range("J74").Value = 32126
$ws.Range("L74").Value = 32126
$ws.Range("N74").Value = -33998
$ws.Range("H77").Value = 32126
$ws.Range("J77").Value = 32126
$ws.Range("L77").Value = 96378
$ws.Range("N77").Value = -105738
